# Add "Preconditions" to Acceptance test Clear Table
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Acceptance test table clear")

$precondition = "Database is accessable`nSystem can verify confirmation codes`nPayment system works as expected`nTable is not canceled or cleared already"

# Rows 2-5 correspond to the four clear-table test cases that were
# missing a Precondition (column D). Use the same style as the
# existing "Expected Result" column (C) on each row so the new cells
# wrap text the same way.
foreach ($r in 2..5) {
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = $precondition
    $dCell.WrapText = $true
}

# Leave the sheet scrolled/selected the way the author left it after typing
# the preconditions in (near the bottom of the used range).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D9").Select()
